$d = $word.ActiveDocument

# --- Paragraph 1: title / subtitle runs ---
$d.Content.Find.Execute(
    "המאמר היומי של מייק ואוראל: 19.01.25", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "המאמר היומי של מייק ואוראל - 18.01.25", 2) | Out-Null

$d.Content.Find.Execute(
    "The Lottery Ticket Hypothesis: Finding Sparse, Trainable Neural Networks", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "MAKING TEXT EMBEDDERS FEW-SHOT LEARNERS", 2) | Out-Null

# --- Paragraph 2 ---
$d.Content.Find.Execute(
    "היפותזת כרטיס הלוטו (Lottery Ticket Hypothesis) אומרת שבתוך רשת נוירונים  צפופה (dense neural nets) המאותחלת בצורה רנדומלית, יש תת-רשת (או ""כרטיס מנצח"") שמאמנים אותה בנפרד, היא יכולה להגיע לביצועים כמו של הרשת המקורית.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "היום להבדיל מהסקירות האחרונות נסקור מאמר מאוד קליל, הלא מערב מתמטיקה כבדה. המאמר מציע שיטה לבניית ייצוג (אמבדינגס) מותאם ללמידה in-context או בקצרה ל-ICL. אזכיר כי ICL היא שיטת בניית פרומפטים כאשר אנו מספקים למודל כמה דוגמאות עבור משימה שאנו מצפים ממנו שיעשה. למשל במשימת גנרוט קוד אנו מספקים למודל (בתוך הפרומפט) כמה דוגמאות שכל אחת מהן היא זוג (שאלה, קוד) במטרה ״להבהיר״ למודל מה אנחנו מצפים ממנו. ד״א למה ICL לפעמים עובד על המשימות שהמודל לא אומן עליהם אינו ברור ב-100% מהווה נושא מחקר די פעיל.",
    2) | Out-Null

# --- Paragraph 3 ---
$d.Content.Find.Execute(
    "נמצא שטכניקת חיתוך(pruning) סטנדרטית מגלה באופן טבעי תת-רשתות כאלה, אשר עבורן מתקיים כי האתחול המחודש תחת אותם hyperparameters, משמר את התוצאות של הרשת המקורית בעלות זולה יותר, כך שהכרטיסים המנצחים הם תת-רשתות אשר ""זכו בהגרלת האתחול"", ובהן המשקלים ההתחלתיים הופכים את האימון לאפקטיבי במיוחד.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "נציין כי המודל בנידון עדיין צריך לגנרט טקסט כלומר יש לנו מודל דקודר (עם מיסוך קוזאלי שדי מפריע לבניית האמבדינג) ונשאלת השאלה איך אנו בונים אמבדינג איתו כמו שאנו רגילים לעשות עם האנקודר. דרך אגב יצאו כמה מאמרים שהציעו שיטות לבניית אמבדינג עם מודלי דקודר כמו LLM2Vec ו-GritML אבל הם אינם מותאמים למקרה שנדון במאמר. כלומר השאלה איך אנו בונים אמבדינג של פרומפט בסגנון ICL כלומר כזה שמכיל כמה דוגמאות פתורות להדגמה. ",
    2) | Out-Null

# --- Paragraph 4 ---
$d.Content.Find.Execute(
    "הרעיון הזה מדגיש את החשיבות של המשקלים ההתחלתיים של הרשת. הכרטיסים המנצחים אינם תת-רשתות אקראיות, אלא כאלה שמתאימות במיוחד בגלל האתחול שלהן. תהליך מציאת התת-רשתות הללו אינו פשוט, כיוון שהוא כרוך בזיהוי החלקים הקריטיים(הנוירונים המשמעותיים) ברשת כבר מההתחלה.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "אז המחברים מצאו לזה פתרון די פשוט. קודם כל הם הוסיפו טוקן EOS בסוף הפרומפט והתכנון הוא שייצוג הטוקן הזה יכיל את האמבדינג של הפרומפט כולו (כמו שנעשה ב-BERT לפני 7 שנים). באופן לא מפתיע המחברים בחרו לעשות זאת עם למידה ניגודית(contrastive learning או CL). מטרה של CL היא לאמן מודל ייצוג כך שהייצוגים של דוגמאות דומות(חיוביות) יהיו קרובות ואילו אלו של דוגמאות לא דומות(שליליות) יהיו רחוקים במרחק האמבדינג. בתור דוגמאות חיוביות המחברים בחרו כאלו עם תשובה נכונה על השאלה בפרומפט ואילו עבור דוגמאות שליליות מופיעות התשובה הלא נכונה. נציין כי הדוגמאות להדגמה בפרומפט נשארות זהות עבור החיוביים והשליליים. ",
    2) | Out-Null

# --- Paragraph 5 ---
$d.Content.Find.Execute(
    "מה זה חיתוך רשת?", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "זהו זה - ככה הם מאמנים מודל אמבדינג על מספר לא גדול של דוגמאות (few-shot) ולפי המאמר התוצאות לא רעות.",
    2) | Out-Null

# --- Remove the five old "pruning deep-dive" paragraphs (old paragraphs 6-10) ---
for ($i = 0; $i -lt 5; $i++) {
    $d.Paragraphs(6).Range.Delete()
}

# --- Final paragraph: swap the arxiv link ---
$d.Content.Find.Execute(
    "https://arxiv.org/pdf/1803.03635", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "https://arxiv.org/abs/2409.15700", 2) | Out-Null
